$wb = $excel.ActiveWorkbook

# --- TransactionTypes: move the cursor/selection to C4 (no change of active sheet) ---
$wsTypes = $wb.Worksheets.Item("TransactionTypes")
$wsTypes.Range("C4").Select()

# --- Wallets: shrink the selection to just C2 and correct the fee value ---
$wsWallets = $wb.Worksheets.Item("Wallets")
$wsWallets.Activate()
$wsWallets.Range("C2").Value = 99
$wsWallets.Range("C2").Select()

# --- WalletTransactions: insert a new "TransactionID" column before the Balance column ---
$wsWT = $wb.Worksheets.Item("WalletTransactions")
$wsWT.Activate()
$wsWT.Columns.Item(4).Insert()
$wsWT.Columns.Item(4).ColumnWidth = $wsWT.Columns.Item(3).ColumnWidth
$wsWT.Range("D1").Value = "TransactionID"

# append the new running-import row
$wsWT.Range("A8").Value = 7
$wsWT.Range("B8").Value = 1
$wsWT.Range("C8").Value = 4
$wsWT.Range("E8").Value = 99
$wsWT.Range("E8").NumberFormat = "0.00"

$wsWT.Range("A1").Select()

# --- Issues: update the running-import quantities ---
$wsIssues = $wb.Worksheets.Item("Issues")
$wsIssues.Range("E3").Value = 4
$wsIssues.Range("E8").Value = 6
$wsIssues.Range("E10").Value = 10

# WalletTransactions is the sheet that should end up active
$wsWT.Activate()
